$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.804.67"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "2.104.15"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +2.31%  "

$ws.Range("E10").Value = "  +2.56%  "

$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").Value = "2.414.92"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").Value = "2.121.67"
$ws.Range("E17").Value = "  +3.83%  "

$ws.Range("D18").Value = "37.649.18"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("E21").Value = "  +1.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.04%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.84%  "

$ws.Range("E37").Value = "  +4.02%  "

$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.90%  "

$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.468.05"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = "  +3.16%  "

$ws.Range("E48").Value = "  +3.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "

$ws.Range("D51").Value = "2.299.93"
$ws.Range("E51").Value = "  +2.31%  "
